# repull data, push all data, mean calculation
# Update the dSF column (F) values for rows where re-pulled data changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 6
    3  = -10
    4  = 5
    11 = -10
    14 = -5
    18 = 4
    20 = -11
    21 = 12
    23 = 3
    24 = -3
    26 = 2
    28 = -5
    29 = -7
    30 = 9
    32 = -3
    40 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
